$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: the 514390ce file reverted from "Ready for handoff"
# back to "In Translation", and the latest HO xliff generation
# timestamp advanced for every row.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Range("G2").Value = "2016-10-24 10:23:05"
$wsOverview.Range("G3").Value = "2016-10-24 10:23:05"
$wsOverview.Range("G4").Value = "2016-10-24 10:23:05"
$wsOverview.Range("G5").Value = "2016-10-24 10:23:05"

# ---------------------------------------------------------------------
# zh-cn sheet: matching status change, priority engine changed from
# "ht" to "mt", and the handoff datetime advanced.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-10-24 10:22:52"
$wsZhCn.Range("H3").Value = "2016-10-24 10:22:52"
$wsZhCn.Range("H4").Value = "2016-10-24 10:22:52"
$wsZhCn.Range("H5").Value = "2016-10-24 10:22:52"

# ---------------------------------------------------------------------
# de-de sheet: same status + priority-engine change, and the handoff
# datetime advanced (shares the Overview sheet's new timestamp).
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-10-24 10:23:05"
$wsDeDe.Range("H3").Value = "2016-10-24 10:23:05"
$wsDeDe.Range("H4").Value = "2016-10-24 10:23:05"
$wsDeDe.Range("H5").Value = "2016-10-24 10:23:05"
